$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (55 cell updates) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 85.666664  # H6: was 96.666664
$ws.Cells.Item(6, 9).Value = 85.666664  # I6: was 96.666664
$ws.Cells.Item(6, 11).Value = 256.999992  # K6: was 289.999992
$ws.Cells.Item(6, 13).Value = -144.999992  # M6: was -177.999992
$ws.Cells.Item(8, 8).Value = 8710.5  # H8: was 10121.667
$ws.Cells.Item(8, 9).Value = 65.25  # I8: was 76.666664
$ws.Cells.Item(8, 10).Value = 26001  # J8: was 20166.666
$ws.Cells.Item(8, 11).Value = 195.75  # K8: was 229.999992
$ws.Cells.Item(8, 12).Value = 78003  # L8: was 60499.99800000001
$ws.Cells.Item(8, 13).Value = -56.75  # M8: was -90.99999199999999
$ws.Cells.Item(8, 14).Value = -78281  # N8: was -60777.99800000001
$ws.Cells.Item(17, 8).Value = 1388.4286  # H17: was 1383.6
$ws.Cells.Item(17, 10).Value = 1388.4286  # J17: was 1383.6
$ws.Cells.Item(17, 12).Value = 4165.2858  # L17: was 4150.799999999999
$ws.Cells.Item(17, 14).Value = -4501.2858  # N17: was -4486.799999999999
$ws.Cells.Item(41, 8).Value = 285.6316  # H41: was 311
$ws.Cells.Item(41, 9).Value = 335  # I41: was 360.625
$ws.Cells.Item(41, 10).Value = 241.2  # J41: was 266.8889
$ws.Cells.Item(41, 11).Value = 335  # K41: was 360.625
$ws.Cells.Item(41, 12).Value = 241.2  # L41: was 266.8889
$ws.Cells.Item(41, 13).Value = 105  # M41: was 79.375
$ws.Cells.Item(41, 14).Value = -1121.2  # N41: was -1146.8889
$ws.Cells.Item(61, 8).Value = 195  # H61: was 590
$ws.Cells.Item(61, 9).Value = 195  # I61: was 590
$ws.Cells.Item(61, 11).Value = 585  # K61: was 1770
$ws.Cells.Item(61, 13).Value = -413  # M61: was -1598
$ws.Cells.Item(75, 8).Value = 34000  # H75: was 32000
$ws.Cells.Item(75, 9).Value = 0  # I75: was 30000
$ws.Cells.Item(75, 10).Value = 34000  # J75: was 32666.666
$ws.Cells.Item(75, 11).Value = 0  # K75: was 30000
$ws.Cells.Item(75, 12).Value = 34000  # L75: was 32666.666
$ws.Cells.Item(75, 13).ClearContents()  # M75: was -29064
$ws.Cells.Item(75, 14).Value = -35872  # N75: was -34538.666
$ws.Cells.Item(78, 8).Value = 34000  # H78: was 32000
$ws.Cells.Item(78, 9).Value = 0  # I78: was 30000
$ws.Cells.Item(78, 10).Value = 34000  # J78: was 32666.666
$ws.Cells.Item(78, 11).Value = 0  # K78: was 90000
$ws.Cells.Item(78, 12).Value = 102000  # L78: was 97999.99800000001
$ws.Cells.Item(78, 13).ClearContents()  # M78: was -85320
$ws.Cells.Item(78, 14).Value = -111360  # N78: was -107359.998
$ws.Cells.Item(113, 8).Value = 3099.5  # H113: was 2857.1428
$ws.Cells.Item(113, 10).Value = 3374.375  # J113: was 3200
$ws.Cells.Item(113, 12).Value = 3374.375  # L113: was 3200
$ws.Cells.Item(113, 14).Value = -9882.375  # N113: was -9708
$ws.Cells.Item(137, 8).Value = 1474.6177  # H137: was 1532.75
$ws.Cells.Item(137, 9).Value = 1245.5927  # I137: was 1301.68
$ws.Cells.Item(137, 11).Value = 3736.7781  # K137: was 3905.04
$ws.Cells.Item(137, 13).Value = -1186.7781  # M137: was -1355.04
$ws.Cells.Item(138, 8).Value = 2153869.5  # H138: was 2329059
$ws.Cells.Item(138, 9).Value = 10527975  # I138: was 11112824
$ws.Cells.Item(138, 10).Value = 3761.3513  # J138: was 3944.7354
$ws.Cells.Item(138, 11).Value = 31583925  # K138: was 33338472
$ws.Cells.Item(138, 12).Value = 11284.0539  # L138: was 11834.2062
$ws.Cells.Item(138, 13).Value = -31578785  # M138: was -33333332
$ws.Cells.Item(138, 14).Value = -21564.0539  # N138: was -22114.2062

# --- Sheet: ARM (54 cell updates) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 12553.646  # H32: was 15485.788
$ws.Cells.Item(32, 9).Value = 13755.667  # I32: was 16616.191
$ws.Cells.Item(32, 10).Value = 3989.25  # J32: was 4860
$ws.Cells.Item(32, 11).Value = 13755.667  # K32: was 16616.191
$ws.Cells.Item(32, 12).Value = 3989.25  # L32: was 4860
$ws.Cells.Item(32, 13).Value = -13468.667  # M32: was -16329.191
$ws.Cells.Item(32, 14).Value = -4563.25  # N32: was -5434
$ws.Cells.Item(63, 8).Value = 3277.7917  # H63: was 3435.318
$ws.Cells.Item(63, 9).Value = 3640.5  # I63: was 3829
$ws.Cells.Item(63, 10).Value = 2770  # J63: was 2866.6667
$ws.Cells.Item(63, 11).Value = 3640.5  # K63: was 3829
$ws.Cells.Item(63, 12).Value = 2770  # L63: was 2866.6667
$ws.Cells.Item(63, 13).Value = -2954.5  # M63: was -3143
$ws.Cells.Item(63, 14).Value = -4142  # N63: was -4238.6667
$ws.Cells.Item(66, 8).Value = 3277.7917  # H66: was 3435.318
$ws.Cells.Item(66, 9).Value = 3640.5  # I66: was 3829
$ws.Cells.Item(66, 10).Value = 2770  # J66: was 2866.6667
$ws.Cells.Item(66, 11).Value = 18202.5  # K66: was 19145
$ws.Cells.Item(66, 12).Value = 13850  # L66: was 14333.3335
$ws.Cells.Item(66, 13).Value = -14770.5  # M66: was -15713
$ws.Cells.Item(66, 14).Value = -20714  # N66: was -21197.3335
$ws.Cells.Item(74, 8).Value = 1109.3  # H74: was 1201.64
$ws.Cells.Item(74, 9).Value = 1006  # I74: was 1066.35
$ws.Cells.Item(74, 10).Value = 1448.7142  # J74: was 1742.8
$ws.Cells.Item(74, 11).Value = 1006  # K74: was 1066.35
$ws.Cells.Item(74, 12).Value = 1448.7142  # L74: was 1742.8
$ws.Cells.Item(74, 13).Value = -132  # M74: was -192.3499999999999
$ws.Cells.Item(74, 14).Value = -3196.7142  # N74: was -3490.8
$ws.Cells.Item(77, 8).Value = 1109.3  # H77: was 1201.64
$ws.Cells.Item(77, 9).Value = 1006  # I77: was 1066.35
$ws.Cells.Item(77, 10).Value = 1448.7142  # J77: was 1742.8
$ws.Cells.Item(77, 11).Value = 5030  # K77: was 5331.75
$ws.Cells.Item(77, 12).Value = 7243.571  # L77: was 8714
$ws.Cells.Item(77, 13).Value = -662  # M77: was -963.75
$ws.Cells.Item(77, 14).Value = -15979.571  # N77: was -17450
$ws.Cells.Item(102, 8).Value = 1965  # H102: was 1225
$ws.Cells.Item(102, 9).Value = 1951.875  # I102: was 1225
$ws.Cells.Item(102, 10).Value = 2000  # J102: was 0
$ws.Cells.Item(102, 11).Value = 1951.875  # K102: was 1225
$ws.Cells.Item(102, 12).Value = 2000  # L102: was 0
$ws.Cells.Item(102, 13).Value = -329.875  # M102: was 397
$ws.Cells.Item(102, 14).Value = -5244  # N102: was None
$ws.Cells.Item(112, 8).Value = 23454  # H112: was 23699.4
$ws.Cells.Item(112, 10).Value = 23454  # J112: was 23699.4
$ws.Cells.Item(112, 12).Value = 23454  # L112: was 23699.4
$ws.Cells.Item(112, 14).Value = -26408  # N112: was -26653.4
$ws.Cells.Item(124, 8).Value = 24000  # H124: was 31464.5
$ws.Cells.Item(124, 10).Value = 24000  # J124: was 31464.5
$ws.Cells.Item(124, 12).Value = 24000  # L124: was 31464.5
$ws.Cells.Item(124, 14).Value = -33820  # N124: was -41284.5
$ws.Cells.Item(125, 8).Value = 72715  # H125: was 72415.836
$ws.Cells.Item(125, 10).Value = 72715  # J125: was 72415.836
$ws.Cells.Item(125, 12).Value = 72715  # L125: was 72415.836
$ws.Cells.Item(125, 14).Value = -82555  # N125: was -82255.836

# --- Sheet: BSM (14 cell updates) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 5001.6665  # H105: was 5010
$ws.Cells.Item(105, 9).Value = 5002  # I105: was 5010
$ws.Cells.Item(105, 10).Value = 5000  # J105: was 0
$ws.Cells.Item(105, 11).Value = 5002  # K105: was 5010
$ws.Cells.Item(105, 12).Value = 5000  # L105: was 0
$ws.Cells.Item(105, 13).Value = -3255  # M105: was -3263
$ws.Cells.Item(105, 14).Value = -8494  # N105: was None
$ws.Cells.Item(134, 8).Value = 3273.2  # H134: was 3276.6
$ws.Cells.Item(134, 9).Value = 3087.375  # I134: was 3056.25
$ws.Cells.Item(134, 10).Value = 3485.5715  # J134: was 3528.4285
$ws.Cells.Item(134, 11).Value = 9262.125  # K134: was 9168.75
$ws.Cells.Item(134, 12).Value = 10456.7145  # L134: was 10585.2855
$ws.Cells.Item(134, 13).Value = -6727.125  # M134: was -6633.75
$ws.Cells.Item(134, 14).Value = -15526.7145  # N134: was -15655.2855

# --- Sheet: CRP (47 cell updates) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(10, 8).Value = 2918.25  # H10: was 12001.4
$ws.Cells.Item(10, 9).Value = 2918.25  # I10: was 3835.6667
$ws.Cells.Item(10, 10).Value = 0  # J10: was 24250
$ws.Cells.Item(10, 11).Value = 2918.25  # K10: was 3835.6667
$ws.Cells.Item(10, 12).Value = 0  # L10: was 24250
$ws.Cells.Item(10, 13).Value = -2779.25  # M10: was -3696.6667
$ws.Cells.Item(10, 14).ClearContents()  # N10: was -24528
$ws.Cells.Item(31, 8).Value = 27030706  # H31: was 32262206
$ws.Cells.Item(31, 9).Value = 100004620  # I31: was 166673260
$ws.Cells.Item(31, 10).Value = 3328.8147  # J31: was 3551.2
$ws.Cells.Item(31, 11).Value = 100004620  # K31: was 166673260
$ws.Cells.Item(31, 12).Value = 3328.8147  # L31: was 3551.2
$ws.Cells.Item(31, 13).Value = -100004325  # M31: was -166672965
$ws.Cells.Item(31, 14).Value = -3918.8147  # N31: was -4141.2
$ws.Cells.Item(34, 8).Value = 27030706  # H34: was 32262206
$ws.Cells.Item(34, 9).Value = 100004620  # I34: was 166673260
$ws.Cells.Item(34, 10).Value = 3328.8147  # J34: was 3551.2
$ws.Cells.Item(34, 11).Value = 100004620  # K34: was 166673260
$ws.Cells.Item(34, 12).Value = 3328.8147  # L34: was 3551.2
$ws.Cells.Item(34, 13).Value = -100004418  # M34: was -166673058
$ws.Cells.Item(34, 14).Value = -3732.8147  # N34: was -3955.2
$ws.Cells.Item(58, 8).Value = 3215.6667  # H58: was 3898
$ws.Cells.Item(58, 9).Value = 3450  # I58: was 4250
$ws.Cells.Item(58, 10).Value = 2747  # J58: was 3194
$ws.Cells.Item(58, 11).Value = 3450  # K58: was 4250
$ws.Cells.Item(58, 12).Value = 2747  # L58: was 3194
$ws.Cells.Item(58, 13).Value = -3247  # M58: was -4047
$ws.Cells.Item(58, 14).Value = -3153  # N58: was -3600
$ws.Cells.Item(99, 8).Value = 3683.7144  # H99: was 3632.96
$ws.Cells.Item(99, 9).Value = 3768.3333  # I99: was 3683.7368
$ws.Cells.Item(99, 11).Value = 3768.3333  # K99: was 3683.7368
$ws.Cells.Item(99, 13).Value = -2270.3333  # M99: was -2185.7368
$ws.Cells.Item(126, 8).Value = 3683.7144  # H126: was 3632.96
$ws.Cells.Item(126, 9).Value = 3768.3333  # I126: was 3683.7368
$ws.Cells.Item(126, 11).Value = 11304.9999  # K126: was 11051.2104
$ws.Cells.Item(126, 13).Value = -8834.999899999999  # M126: was -8581.2104
$ws.Cells.Item(136, 8).Value = 3215.6667  # H136: was 3898
$ws.Cells.Item(136, 9).Value = 3450  # I136: was 4250
$ws.Cells.Item(136, 10).Value = 2747  # J136: was 3194
$ws.Cells.Item(136, 11).Value = 10350  # K136: was 12750
$ws.Cells.Item(136, 12).Value = 8241  # L136: was 9582
$ws.Cells.Item(136, 13).Value = -7800  # M136: was -10200
$ws.Cells.Item(136, 14).Value = -13341  # N136: was -14682
$ws.Cells.Item(141, 8).Value = 45163  # H141: was 47744.25
$ws.Cells.Item(141, 10).Value = 50326  # J141: was 50325.668
$ws.Cells.Item(141, 12).Value = 50326  # L141: was 50325.668
$ws.Cells.Item(141, 14).Value = -60686  # N141: was -60685.668

# --- Sheet: CUL (33 cell updates) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(15, 8).Value = 634.125  # H15: was 205.2
$ws.Cells.Item(15, 9).Value = 11.5  # I15: was 75
$ws.Cells.Item(15, 10).Value = 841.6667  # J15: was 400.5
$ws.Cells.Item(15, 11).Value = 34.5  # K15: was 225
$ws.Cells.Item(15, 12).Value = 2525.0001  # L15: was 1201.5
$ws.Cells.Item(15, 13).Value = 105.5  # M15: was -85
$ws.Cells.Item(15, 14).Value = -2805.0001  # N15: was -1481.5
$ws.Cells.Item(68, 8).Value = 1425.3  # H68: was 1455.3
$ws.Cells.Item(68, 10).Value = 1437.875  # J68: was 1475.375
$ws.Cells.Item(68, 12).Value = 4313.625  # L68: was 4426.125
$ws.Cells.Item(68, 14).Value = -5935.625  # N68: was -6048.125
$ws.Cells.Item(71, 8).Value = 1425.3  # H71: was 1455.3
$ws.Cells.Item(71, 10).Value = 1437.875  # J71: was 1475.375
$ws.Cells.Item(71, 12).Value = 12940.875  # L71: was 13278.375
$ws.Cells.Item(71, 14).Value = -21052.875  # N71: was -21390.375
$ws.Cells.Item(121, 8).Value = 1215.8889  # H121: was 72594.71000000001
$ws.Cells.Item(121, 9).Value = 1195  # I121: was 1326.6666
$ws.Cells.Item(121, 10).Value = 1232.6  # J121: was 126045.75
$ws.Cells.Item(121, 11).Value = 3585  # K121: was 3979.9998
$ws.Cells.Item(121, 12).Value = 3697.8  # L121: was 378137.25
$ws.Cells.Item(121, 13).Value = -2275  # M121: was -2669.9998
$ws.Cells.Item(121, 14).Value = -6317.799999999999  # N121: was -380757.25
$ws.Cells.Item(129, 8).Value = 1853039.8  # H129: was 1786935.1
$ws.Cells.Item(129, 9).Value = 496  # I129: was 517.7778
$ws.Cells.Item(129, 10).Value = 2942771.2  # J129: was 2633132.8
$ws.Cells.Item(129, 11).Value = 1488  # K129: was 1553.3334
$ws.Cells.Item(129, 12).Value = 8828313.600000001  # L129: was 7899398.399999999
$ws.Cells.Item(129, 13).Value = 3512  # M129: was 3446.6666
$ws.Cells.Item(129, 14).Value = -8838313.600000001  # N129: was -7909398.399999999
$ws.Cells.Item(131, 8).Value = 17244470  # H131: was 32263008
$ws.Cells.Item(131, 10).Value = 20409716  # J131: was 45456820
$ws.Cells.Item(131, 12).Value = 61229148  # L131: was 136370460
$ws.Cells.Item(131, 14).Value = -61239228  # N131: was -136380540

# --- Sheet: GSM (16 cell updates) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(114, 8).Value = 0  # H114: was 19900
$ws.Cells.Item(114, 10).Value = 0  # J114: was 19900
$ws.Cells.Item(114, 12).Value = 0  # L114: was 19900
$ws.Cells.Item(114, 14).ClearContents()  # N114: was -28578
$ws.Cells.Item(123, 8).Value = 13666.5  # H123: was 13103.667
$ws.Cells.Item(123, 10).Value = 13666.5  # J123: was 13103.667
$ws.Cells.Item(123, 12).Value = 13666.5  # L123: was 13103.667
$ws.Cells.Item(123, 14).Value = -18566.5  # N123: was -18003.667
$ws.Cells.Item(126, 8).Value = 2778.1177  # H126: was 2393.2273
$ws.Cells.Item(126, 9).Value = 2016.909  # I126: was 1725.5625
$ws.Cells.Item(126, 11).Value = 6050.727000000001  # K126: was 5176.6875
$ws.Cells.Item(126, 13).Value = -3580.727000000001  # M126: was -2706.6875
$ws.Cells.Item(132, 8).Value = 2207.0417  # H132: was 2554
$ws.Cells.Item(132, 9).Value = 1619.579  # I132: was 1880.6428
$ws.Cells.Item(132, 11).Value = 4858.737  # K132: was 5641.928400000001
$ws.Cells.Item(132, 13).Value = -2328.737  # M132: was -3111.928400000001

# --- Sheet: LTW (25 cell updates) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 5014.2856  # H16: was 4271.5713
$ws.Cells.Item(16, 9).Value = 5720  # I16: was 4900
$ws.Cells.Item(16, 10).Value = 3250  # J16: was 501
$ws.Cells.Item(16, 11).Value = 5720  # K16: was 4900
$ws.Cells.Item(16, 12).Value = 3250  # L16: was 501
$ws.Cells.Item(16, 13).Value = -5550  # M16: was -4730
$ws.Cells.Item(16, 14).Value = -3590  # N16: was -841
$ws.Cells.Item(82, 8).Value = 2256.56  # H82: was 2096.6191
$ws.Cells.Item(82, 9).Value = 1813.8  # I82: was 1938.3077
$ws.Cells.Item(82, 10).Value = 2920.7  # J82: was 2353.875
$ws.Cells.Item(82, 11).Value = 1813.8  # K82: was 1938.3077
$ws.Cells.Item(82, 12).Value = 2920.7  # L82: was 2353.875
$ws.Cells.Item(82, 13).Value = -1452.8  # M82: was -1577.3077
$ws.Cells.Item(82, 14).Value = -3642.7  # N82: was -3075.875
$ws.Cells.Item(85, 8).Value = 2256.56  # H85: was 2096.6191
$ws.Cells.Item(85, 9).Value = 1813.8  # I85: was 1938.3077
$ws.Cells.Item(85, 10).Value = 2920.7  # J85: was 2353.875
$ws.Cells.Item(85, 11).Value = 1813.8  # K85: was 1938.3077
$ws.Cells.Item(85, 12).Value = 2920.7  # L85: was 2353.875
$ws.Cells.Item(85, 13).Value = -565.8  # M85: was -690.3077000000001
$ws.Cells.Item(85, 14).Value = -5416.7  # N85: was -4849.875
$ws.Cells.Item(100, 8).Value = 5803.5  # H100: was 8225.375
$ws.Cells.Item(100, 9).Value = 6793.5557  # I100: was 11460.6
$ws.Cells.Item(100, 11).Value = 6793.5557  # K100: was 11460.6
$ws.Cells.Item(100, 13).Value = -6252.5557  # M100: was -10919.6

# --- Sheet: WVR (22 cell updates) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 1096.6  # H96: was 839.55554
$ws.Cells.Item(96, 9).Value = 901.5  # I96: was 536.5714
$ws.Cells.Item(96, 10).Value = 1226.6666  # J96: was 1900
$ws.Cells.Item(96, 11).Value = 901.5  # K96: was 536.5714
$ws.Cells.Item(96, 12).Value = 1226.6666  # L96: was 1900
$ws.Cells.Item(96, 13).Value = 471.5  # M96: was 836.4286
$ws.Cells.Item(96, 14).Value = -3972.6666  # N96: was -4646
$ws.Cells.Item(122, 8).Value = 13160067  # H122: was 22730482
$ws.Cells.Item(122, 9).Value = 17859082  # I122: was 31253290
$ws.Cells.Item(122, 10).Value = 2824  # J122: was 2995
$ws.Cells.Item(122, 11).Value = 53577246  # K122: was 93759870
$ws.Cells.Item(122, 12).Value = 8472  # L122: was 8985
$ws.Cells.Item(122, 13).Value = -53574796  # M122: was -93757420
$ws.Cells.Item(122, 14).Value = -13372  # N122: was -13885
$ws.Cells.Item(132, 8).Value = 1371.4642  # H132: was 1154.0294
$ws.Cells.Item(132, 9).Value = 1018.1818  # I132: was 829.8570999999999
$ws.Cells.Item(132, 11).Value = 3054.5454  # K132: was 2489.5713
$ws.Cells.Item(132, 13).Value = -524.5454  # M132: was 40.42870000000039
$ws.Cells.Item(136, 8).Value = 1128.25  # H136: was 1140.1305
$ws.Cells.Item(136, 9).Value = 1112.0869  # I136: was 1123.7727
$ws.Cells.Item(136, 11).Value = 3336.2607  # K136: was 3371.3181
$ws.Cells.Item(136, 13).Value = -786.2606999999998  # M136: was -821.3181
